$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Status" column (B) is filled in with the patient counts that were
# collected for every remaining country, marking each row as completed
# (green fill), and the orange "in progress" marker on Kazakhstan (B7)
# is replaced with a normal completed count as well.

$green = [System.Drawing.ColorTranslator]::ToOle([System.Drawing.Color]::FromArgb(0, 176, 80))

$values = @{
    2  = 345   # Azerbaijan
    4  = 67    # China
    6  = 15    # India
    7  = 974   # Kazakhstan
    8  = 37    # Kyrgyzstan
    10 = 756   # Nigeria
    11 = 465   # Romania
    12 = 2     # Senegal
    13 = 101   # South Africa
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $values[$row]
    $cell.Interior.Color = $green
}

# Update the total so it also includes Azerbaijan (row 2).
$ws.Range("B15").Formula = "=SUM(B2:B14)"

# Update the saved cursor/selection position.
$ws.Range("E10").Select()
